$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update Q_cool value (B3). B5 formula (B3/B4) will recalculate automatically.
$ws.Range("B3").Value = 110905.980975

$wb.Application.Calculate()
